# Generate Report for Handoff
# - Update status text from "Handed back: in sync with en-US" to "Ready for handoff"
# - Refresh the associated timestamps
# - Shrink the now-shorter status columns to fit the new text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Timestamps refreshed to reflect the new handoff generation
$overview.Range("G2").Value = "2016-08-24 13:04:15"
$dede.Range("H2").Value = "2016-08-24 13:04:15"
$zhcn.Range("H2").Value = "2016-08-24 13:04:03"

# Column widths shrink now that the status text is shorter.
# (ColumnWidth is quantized to whole pixels by Excel, so we pick the input
# that lands on the closest achievable stored width to the target 17.2159881591797.)
$overview.Range("E1").ColumnWidth = 16.3333333333333
$overview.Range("F1").ColumnWidth = 16.3333333333333
$zhcn.Range("C1").ColumnWidth = 16.3333333333333
$dede.Range("C1").ColumnWidth = 16.3333333333333
